$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Chocolate"
$ws.Range("B1").Value = 13243
$ws.Range("C1").Value = 30.5
$ws.Range("D1").Value = 23
$ws.Range("E1").Value = 1500
$ws.Range("F1").Value = 45251
$ws.Range("F1").NumberFormat = "mm-dd-yy"

$ws.Range("A2").Value = "Crisp Chips"
$ws.Range("B2").Value = 2346
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 18.5
$ws.Range("E2").Value = 300
$ws.Range("F2").Value = 45291

$ws.Range("F1").Copy()
$ws.Range("F2").PasteSpecial(-4122)

$ws.Range("D2").Select()
